$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers (H1:J1) for the brute-force clustering columns ---
$ws.Range("H1").Value = "ClusterSize(Brute-force)"
$ws.Range("I1").Value = "MinimumPercentIdentity(Brute-force)"
$ws.Range("J1").Value = "Average(Brute-force)"

# Copy the existing header formatting (fill colour) onto the new headers
$ws.Range("A1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New data for rows 2-51 (columns H, I, J) ---
$data = @(
    @(86, 72, 85.7751),
    @(82, 68, 85.0425),
    @(71, 83, 93.1505),
    @(58, 78, 94.9625),
    @(58, 62, 77.8264),
    @(55, 59, 75.9313),
    @(54, 65, 95.2229),
    @(53, 72, 88.1255),
    @(53, 87, 95.4935),
    @(48, 65, 77.9441),
    @(48, 69, 77.3777),
    @(45, 67, 74.3242),
    @(45, 67, 88.6232),
    @(45, 78, 92.5091),
    @(44, 60, 75.6554),
    @(42, 65, 87.3008),
    @(41, 65, 74.5793),
    @(39, 68, 78.2645),
    @(39, 68, 76.9838),
    @(39, 61, 76.9892),
    @(38, 69, 92.7795),
    @(38, 69, 92.1863),
    @(38, 65, 78.0612),
    @(37, 78, 88.5736),
    @(36, 69, 78.9667),
    @(36, 81, 90.9635),
    @(36, 63, 83.5698),
    @(34, 79, 92.1836),
    @(34, 57, 78.123),
    @(33, 61, 78.7235),
    @(32, 65, 80.6371),
    @(32, 58, 84.8528),
    @(32, 72, 82.0968),
    @(32, 67, 81.2319),
    @(32, 67, 80.2883),
    @(32, 75, 87.0464),
    @(31, 75, 88.2645),
    @(31, 68, 81.6086),
    @(31, 62, 80.8645),
    @(31, 66, 86.0323),
    @(31, 67, 79.6989),
    @(30, 66, 84.4345),
    @(30, 69, 85.092),
    @(30, 71, 87.8046),
    @(29, 67, 77.8103),
    @(29, 69, 83.4039),
    @(28, 66, 77.2725),
    @(28, 60, 77.1693),
    @(28, 64, 76.418),
    @(28, 68, 78.6005)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $vals = $data[$idx]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
    $ws.Cells.Item($row, 10).Value = $vals[2]
}

# --- Resize columns to fit new content, matching the post-edit layout ---
$ws.Range("A1:J51").EntireColumn.AutoFit() | Out-Null

# --- Update selection to match the post-edit state ---
$ws.Range("A2:J51").Select() | Out-Null
